$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": row 3 (b.md) is now ready for handoff.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-33-13 04:33:01"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": row 3 (b.md) got a new handoff file + status + datetime.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-13 04:32:57"

# Rebuild the hyperlinks collection so the D3 hyperlink's display text
# reflects the new handoff file while every other hyperlink (and the D3
# link target itself) stays exactly as it was.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09d4d05d9e0d5bc56eb4ca45c187e0db11553c84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6038559a70152a0ca97a9be52ab6b165caee6048/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f85cc946da0eda73827198d9d8ee3432ee8fe85/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09d4d05d9e0d5bc56eb4ca45c187e0db11553c84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6038559a70152a0ca97a9be52ab6b165caee6048/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f85cc946da0eda73827198d9d8ee3432ee8fe85/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de": row 3 (b.md) got a new handoff file + status + datetime.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-13 04:33:01"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/a.md", [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b59b8fb3016e9c272087c4b349d20bd6b2f5a4ab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/69df9c36eb548322289f47b8317ac6b64c3c97d9/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/78ea9d62edd977d01436aeb6fd68e93a491791bb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/65e90de971cf4d9adc07e5045d8795d5b633077e/e2e/b.md", [Type]::Missing, [Type]::Missing, ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b59b8fb3016e9c272087c4b349d20bd6b2f5a4ab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/69df9c36eb548322289f47b8317ac6b64c3c97d9/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/78ea9d62edd977d01436aeb6fd68e93a491791bb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

Write-Output "Report regenerated for handoff."
